$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.714.78'
$ws.Range("E2").Value = '  +3.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.010.91'
$ws.Range("E3").Value = '  +3.25%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.39'
$ws.Range("E5").Value = '  +3.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.51'
$ws.Range("E6").Value = '  +8.65%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +2.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.997.58'
$ws.Range("E9").Value = '  +3.00%  '
$ws.Range("E10").Value = '  +6.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.30'
$ws.Range("E11").Value = '  +12.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.455'
$ws.Range("E12").Value = '  +2.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").Value = '  +6.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.06'
$ws.Range("E14").Value = '  +4.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.122'
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.505.61'
$ws.Range("E16").Value = '  +3.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.24'
$ws.Range("E17").Value = '  +5.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.006.80'
$ws.Range("E18").Value = '  +3.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '59.695.25'
$ws.Range("E19").Value = '  +3.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '439.48'
$ws.Range("E20").Value = '  +5.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.66'
$ws.Range("E21").Value = '  +3.34%  '
$ws.Range("E22").Value = '  +4.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.57'
$ws.Range("E23").Value = '  +3.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.12'
$ws.Range("E24").Value = '  +2.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.62'
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  +11.80%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.56'
$ws.Range("E29").Value = '  +3.92%  '
$ws.Range("E30").Value = '  +6.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.29'
$ws.Range("E31").Value = '  +5.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.92'
$ws.Range("E32").Value = '  +3.12%  '
$ws.Range("E33").Value = '  +10.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0785'
$ws.Range("E34").Value = '  +16.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +7.92%  '
$ws.Range("E36").Value = '  +5.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.10'
$ws.Range("E37").Value = '  +2.71%  '
$ws.Range("E38").Value = '  +2.11%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.60'
$ws.Range("E39").Value = '  -1.19%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.80'
$ws.Range("E40").Value = '  +10.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '404.26'
$ws.Range("E41").Value = '  +7.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0356'
$ws.Range("E42").Value = '  +4.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.771.81'
$ws.Range("E43").Value = '  +3.20%  '
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("E45").Value = '  +7.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '123.48'
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("E48").Value = '  +5.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.53'
$ws.Range("E49").Value = '  +21.74%  '
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.83'
$ws.Range("E51").Value = '  +4.24%  '
